# V1.5.3 - add new numeric metrics.
# Adds a new round ("Cypress Point") to the "Score Cards" table and a
# matching new course-par line to the "Course Pars" table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Score Cards" sheet / Table1 -> append row 16 (Cypress Point round)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Score Cards")

# Clone the formatting of the previous last row (row 15) onto the new
# row (row 16) so that date/number styles carry over correctly.
$ws1.Range("A15:W15").Copy()
$ws1.Range("A16:W16").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A16").Value = "Cypress Point"
$ws1.Range("B16").Value = 45899

$row16Cols  = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
$row16Vals  = @(6,6,3,7,5,6,5,6,4,5,3,4,5,5,4,4,5,6,5,37,6)
for ($i = 0; $i -lt $row16Cols.Length; $i++) {
    $ws1.Range("$($row16Cols[$i])16").Value = $row16Vals[$i]
}

# Grow Table1 so the new row becomes part of the table.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:W16"))

$ws1.Activate() | Out-Null
$ws1.Range("A17").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "Course Pars" sheet / Table24 -> append row 4 (Cypress Point pars)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Course Pars")

# Clone the formatting of the previous last row (row 3) onto the new
# row (row 4).
$ws3.Range("A3:U3").Copy()
$ws3.Range("A4:U4").PasteSpecial(-4122)  # xlPasteFormats

$ws3.Range("A4").Value = "Cypress Point"

$row4Cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")
$row4Vals = @(4,5,3,5,3,4,4,5,3,4,3,4,4,5,4,4,3,5,69,118)
for ($i = 0; $i -lt $row4Cols.Length; $i++) {
    $ws3.Range("$($row4Cols[$i])4").Value = $row4Vals[$i]
}

# Grow Table24 so the new row becomes part of the table.
$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:U4"))

$ws3.Activate() | Out-Null
$ws3.Range("A5").Select() | Out-Null

# Leave "Score Cards" as the active / selected sheet, matching the
# original workbook's tab selection.
$ws1.Activate() | Out-Null
